# Append two new trading-log rows (102 and 103) to Sheet1, matching the
# existing log format used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 102: TRADING_ATTEMPT for AAVE
$ws.Cells.Item(102, 1).Value = "2025-11-09T01:42:41.703833"
$ws.Cells.Item(102, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(102, 3).Value = "AAVE"
$ws.Cells.Item(102, 4).Value = "UNKNOWN"
$ws.Cells.Item(102, 5).Value = 200.8622597986005
$ws.Cells.Item(102, 6).Value = ""
$ws.Cells.Item(102, 7).Value = ""
$ws.Cells.Item(102, 8).Value = ""
$ws.Cells.Item(102, 9).Value = ""
$ws.Cells.Item(102, 10).Value = ""
$ws.Cells.Item(102, 11).Value = "ATTEMPT"
$ws.Cells.Item(102, 12).Value = "Attempting trade 1/1"

# Row 103: POSITION_FAILED for AAVE
$ws.Cells.Item(103, 1).Value = "2025-11-09T01:42:43.221998"
$ws.Cells.Item(103, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(103, 3).Value = "AAVE"
$ws.Cells.Item(103, 4).Value = "UNKNOWN"
$ws.Cells.Item(103, 5).Value = ""
$ws.Cells.Item(103, 6).Value = ""
$ws.Cells.Item(103, 7).Value = ""
$ws.Cells.Item(103, 8).Value = ""
$ws.Cells.Item(103, 9).Value = ""
$ws.Cells.Item(103, 10).Value = ""
$ws.Cells.Item(103, 11).Value = "FAILED"
$ws.Cells.Item(103, 12).Value = "Trade execution failed for trade 1"
